$p = $ppt.ActivePresentation

# Remove the duplicate "COMMUNITY PROTECTION" slide (position 6 in the deck) --
# it was a near-duplicate (subset) of the following slide and got cleaned up.
$p.Slides.Item(6).Delete()

# On the "HANDLING CHEMICALS" slide (now slide 12 after the deletion above),
# split the sentence "Treat all chemicals as a safety" so it reads
# "Treat all chemicals as a safety concern".
$s = $p.Slides.Item(12)
$shp = $s.Shapes.Item(2)
$para = $shp.TextFrame.TextRange.Paragraphs(2, 1)
$word = $para.Characters(26, 6)
$word.Text = "safety concern"
